# Add "aansluiting" breakdown rows to the levels table and fix row 12
# (which previously only had column A filled in - "dubbel" / duplicate bug).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 12-19 (columns A-G)
# A = VAR_Formal_variable, B = VAR_Simple_variable, C = VAR_Level_NL,
# D = VAR_Level_label_NL, E = VAR_Level_label_NL_description,
# F = VAR_Level_order, G = VAR_Breakdown_tf
$rows = @(
    @("aansluiting", "aansluiting", "Direct",        "Direct",        "studenten die direct na hun vooropleiding instromen", 1, "TRUE"),
    @("aansluiting", "aansluiting", "Tussenjaar",     "Tussenjaar",    "studenten met een of meer tussenjaren", 2, "TRUE"),
    @("aansluiting", "aansluiting", "Switch intern",  "Switch intern", "interne switchers", 3, "TRUE"),
    @("aansluiting", "aansluiting", "Switch extern",  "Switch extern", "externe switchers", 4, "TRUE"),
    @("aansluiting", "aansluiting", "2e Studie",      "2e Studie",     "studenten die twee of meer studies volgen", 5, "TRUE"),
    @("aansluiting", "aansluiting", "Na CD",          "Na CD",         "studenten die instromen met een 21+ toets of Colloquium Doctum", 6, "TRUE"),
    @("aansluiting", "aansluiting", "Overig",         "Overig",        "studenten met een overige aansluiting", 7, "TRUE"),
    @("aansluiting", "aansluiting", "Onbekend",       "Onbekend",      "studenten met een onbekende aansluiting", 8, "TRUE")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    # Column G holds the literal text "TRUE" (matches existing rows which
    # store it as shared-string text, not a boolean) - use a leading
    # apostrophe to force text entry, then reset the style so no stray
    # quote-prefix formatting lingers on the cell.
    $ws.Cells.Item($r, 7).Value = "'" + $data[6]
    $ws.Cells.Item($r, 7).Style = "Normal"
}

$ws.Range("E14").Select()
